$wb = $excel.ActiveWorkbook

# --- Productdata sheet: scale InventoryCosts (D), BackorderCosts (F), LostSale (I) by 1/2500 ---
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("D2").Value = 0.0016
$ws.Range("F2").Value = 0.016
$ws.Range("I2").Value = 0.16

$ws.Range("D3").Value = 0.0028
$ws.Range("F3").Value = 0.028
$ws.Range("I3").Value = 0.28

$ws.Range("D4").Value = 0.0024
$ws.Range("F4").Value = 0.024
$ws.Range("I4").Value = 0.24

$ws.Range("D5").Value = 0.0012
$ws.Range("F5").Value = 0.012
$ws.Range("I5").Value = 0.12

$ws.Range("D6").Value = 0.0012
$ws.Range("F6").Value = 0.012
$ws.Range("I6").Value = 0.12

$ws.Range("D7").Value = 0.0012
$ws.Range("F7").Value = 0.012
$ws.Range("I7").Value = 0.12

$ws.Range("D8").Value = 0.0008
$ws.Range("F8").Value = 0.008
$ws.Range("I8").Value = 0.08

$ws.Range("D9").Value = 0.0004
$ws.Range("F9").Value = 0.004
$ws.Range("I9").Value = 0.04

$ws.Range("D10").Value = 0.0004
$ws.Range("F10").Value = 0.004
$ws.Range("I10").Value = 0.04

$ws.Range("D11").Value = 0.0004
$ws.Range("F11").Value = 0.004
$ws.Range("I11").Value = 0.04

# --- ForcastedStandardDeviation sheet: zero out rows 9-11, columns B:E ---
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")

$ws2.Range("B9:E9").Value = 0
$ws2.Range("B10:E10").Value = 0
$ws2.Range("B11:E11").Value = 0
